$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.541.85'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.02%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.431.35'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.64%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '503.65'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.60'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.99%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.549'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.57%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.443.73'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.23%  '
$ws.Range("E10").Value = '  -0.23%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0950'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.23%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.17'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.76%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.329'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.05%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.865.50'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.54%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '57.474.75'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.70'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.56%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000132'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.05%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.435.07'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.57%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.41'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '313.97'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.08%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.09'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.17%  '
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.66'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.66%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.19'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.80%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.406'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.59%  '
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.159'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.47%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.20'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.71%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '169.89'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.93%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0719'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.26%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.16'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.21%  '
$ws.Range("E32").Value = '  -3.27%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.12'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.18%  '
$ws.Range("E34").Value = '  -0.03%  '
$ws.Range("E35").Value = '  -0.10%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.70'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.43%  '
$ws.Range("E37").Value = '  -5.29%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.89'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.42%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.45'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.12%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.44'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.09%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.750'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.53%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '269.97'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.50%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.35'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.23%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.89'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.98%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.578'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.92%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0909'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '119.01'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.58%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0483'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.83%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '17.07'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.49%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0208'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.77%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '16.52'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.43%  '
